$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its text formatting so values like
# "66.648.95" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '66.648.95'
$ws.Cells.Item(2, 5).Value = '  +0.13%  '
$ws.Cells.Item(3, 4).Value = '3.466.10'
$ws.Cells.Item(3, 5).Value = '  -1.04%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).Value = '585.39'
$ws.Cells.Item(5, 5).Value = '  +0.11%  '
$ws.Cells.Item(6, 4).Value = '177.96'
$ws.Cells.Item(6, 5).Value = '  +1.64%  '
$ws.Cells.Item(7, 4).Value = '0.628'
$ws.Cells.Item(7, 5).Value = '  +5.50%  '
$ws.Cells.Item(8, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 4).Value = '3.461.78'
$ws.Cells.Item(9, 5).Value = '  -0.99%  '
$ws.Cells.Item(10, 4).Value = '0.134'
$ws.Cells.Item(10, 5).Value = '  +0.87%  '
$ws.Cells.Item(11, 4).Value = '7.00'
$ws.Cells.Item(11, 5).Value = '  +1.43%  '
$ws.Cells.Item(12, 4).Value = '0.420'
$ws.Cells.Item(12, 5).Value = '  -0.27%  '
$ws.Cells.Item(13, 4).Value = '4.058.84'
$ws.Cells.Item(13, 5).Value = '  -1.18%  '
$ws.Cells.Item(14, 4).Value = '0.135'
$ws.Cells.Item(14, 5).Value = '  +1.48%  '
$ws.Cells.Item(15, 4).Value = '30.13'
$ws.Cells.Item(15, 5).Value = '  -0.52%  '
$ws.Cells.Item(16, 4).Value = '66.520.74'
$ws.Cells.Item(16, 5).Value = '  -0.05%  '
$ws.Cells.Item(17, 4).Value = '0.0000174'
$ws.Cells.Item(17, 5).Value = '  +1.12%  '
$ws.Cells.Item(18, 4).Value = '3.451.52'
$ws.Cells.Item(18, 5).Value = '  -1.28%  '
$ws.Cells.Item(19, 4).Value = '5.99'
$ws.Cells.Item(19, 5).Value = '  -0.26%  '
$ws.Cells.Item(20, 4).Value = '13.92'
$ws.Cells.Item(20, 5).Value = '  +0.59%  '
$ws.Cells.Item(21, 4).Value = '372.21'
$ws.Cells.Item(21, 5).Value = '  -2.09%  '
$ws.Cells.Item(22, 4).Value = '7.69'
$ws.Cells.Item(22, 5).Value = '  -2.11%  '
$ws.Cells.Item(23, 4).Value = '73.55'
$ws.Cells.Item(23, 5).Value = '  +1.99%  '
$ws.Cells.Item(24, 4).Value = '1.00'
$ws.Cells.Item(24, 5).Value = '  +0.07%  '
$ws.Cells.Item(25, 4).Value = '0.538'
$ws.Cells.Item(25, 5).Value = '  -1.88%  '
$ws.Cells.Item(26, 4).Value = '0.0000127'
$ws.Cells.Item(26, 5).Value = '  +5.51%  '
$ws.Cells.Item(27, 4).Value = '10.01'
$ws.Cells.Item(27, 5).Value = '  +1.58%  '
$ws.Cells.Item(28, 4).Value = '0.178'
$ws.Cells.Item(28, 5).Value = '  +2.46%  '
$ws.Cells.Item(29, 5).Value = '  +0.02%  '
$ws.Cells.Item(30, 4).Value = '5.99'
$ws.Cells.Item(30, 5).Value = '  +1.77%  '
$ws.Cells.Item(31, 4).Value = '2.01'
$ws.Cells.Item(31, 5).Value = '  +0.11%  '
$ws.Cells.Item(32, 4).Value = '23.70'
$ws.Cells.Item(32, 5).Value = '  -2.88%  '
$ws.Cells.Item(33, 5).Value = '  +0.00%  '
$ws.Cells.Item(34, 4).Value = '7.09'
$ws.Cells.Item(34, 5).Value = '  -1.94%  '
$ws.Cells.Item(35, 4).Value = '1.28'
$ws.Cells.Item(35, 5).Value = '  -3.23%  '
$ws.Cells.Item(36, 4).Value = '1.58'
$ws.Cells.Item(36, 5).Value = '  +0.39%  '
$ws.Cells.Item(37, 4).Value = '162.42'
$ws.Cells.Item(37, 5).Value = '  +1.54%  '
$ws.Cells.Item(38, 4).Value = '0.886'
$ws.Cells.Item(38, 5).Value = '  -0.63%  '
$ws.Cells.Item(39, 4).Value = '27.96'
$ws.Cells.Item(39, 5).Value = '  -5.45%  '
$ws.Cells.Item(40, 4).Value = '1.82'
$ws.Cells.Item(40, 5).Value = '  +2.07%  '
$ws.Cells.Item(41, 4).Value = '4.53'
$ws.Cells.Item(41, 5).Value = '  +0.19%  '
$ws.Cells.Item(42, 4).Value = '2.776.66'
$ws.Cells.Item(42, 5).Value = '  +3.67%  '
$ws.Cells.Item(43, 4).Value = '2.58'
$ws.Cells.Item(43, 5).Value = '  +1.77%  '
$ws.Cells.Item(44, 4).Value = '6.48'
$ws.Cells.Item(44, 5).Value = '  +0.70%  '
$ws.Cells.Item(45, 4).Value = '0.0697'
$ws.Cells.Item(45, 5).Value = '  +0.14%  '
$ws.Cells.Item(46, 4).Value = '25.52'
$ws.Cells.Item(46, 5).Value = '  +4.59%  '
$ws.Cells.Item(47, 4).Value = '341.56'
$ws.Cells.Item(47, 5).Value = '  +8.34%  '
$ws.Cells.Item(48, 4).Value = '40.11'
$ws.Cells.Item(48, 5).Value = '  -1.56%  '
$ws.Cells.Item(49, 4).Value = '0.0290'
$ws.Cells.Item(49, 5).Value = '  -0.45%  '
$ws.Cells.Item(50, 5).Value = '  +2.77%  '
$ws.Cells.Item(51, 2).Value = 'Arweave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(51, 4).Value = '31.84'
$ws.Cells.Item(51, 5).Value = '  +3.72%  '
